$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Helper: write a "date-looking" text value into a cell while keeping it
# as literal text (not an Excel date serial number), and without leaving
# a stray custom number-format style behind on the cell.
# ---------------------------------------------------------------------
function Set-TextValue {
    param($cell, $value)
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

# =======================================================================
# Sheet "展览" (Exhibition) — "想去人数" (want-to-go count) bumps
# =======================================================================
$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Cells.Item(2, 6).Value = 5307
$wsExpo.Cells.Item(4, 6).Value = 634
$wsExpo.Cells.Item(6, 6).Value = 806
$wsExpo.Cells.Item(7, 6).Value = 304
$wsExpo.Cells.Item(8, 6).Value = 12

# =======================================================================
# Sheet "演出" (Performance) — append a new event row (row 4)
# =======================================================================
$wsShow = $wb.Worksheets.Item("演出")

$wsShow.Cells.Item(4, 1).Value = 3
$wsShow.Range("A2").Copy()
$wsShow.Range("A4").PasteSpecial(-4122)

Set-TextValue $wsShow.Cells.Item(4, 2) "2024-11-09"
$wsShow.Cells.Item(4, 3).Value = "合肥·一生必听的钢琴曲—“从巴赫 · 莫扎特到肖邦 · 李斯特”钢琴圣手谭小棠独奏音乐会"
$wsShow.Cells.Item(4, 4).Value = "徽州大道辅路与祁门路辅路交叉口北120米 包河凤凰剧院"
$wsShow.Cells.Item(4, 5).Value = "2024.11.09 19:30-11.09 21:00"
$wsShow.Cells.Item(4, 6).Value = 0
$wsShow.Cells.Item(4, 7).Value = 56
$wsShow.Cells.Item(4, 8).Value = "https://show.bilibili.com/platform/detail.html?id=90593"
$wsShow.Cells.Item(4, 9).Value = "//i2.hdslb.com/bfs/openplatform/202408/SYfLxnO21723442234232.jpeg"

# =======================================================================
# Sheet "本地生活" (Local life) — unchanged
# =======================================================================

# =======================================================================
# Sheet "全部类型" (All types) — same "想去人数" bumps as 展览, plus the
# same new event appended as a new row (row 11)
# =======================================================================
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Cells.Item(2, 6).Value = 5307
$wsAll.Cells.Item(4, 6).Value = 634
$wsAll.Cells.Item(6, 6).Value = 806
$wsAll.Cells.Item(8, 6).Value = 304
$wsAll.Cells.Item(9, 6).Value = 12

$wsAll.Cells.Item(11, 1).Value = 10
$wsAll.Range("A2").Copy()
$wsAll.Range("A11").PasteSpecial(-4122)

Set-TextValue $wsAll.Cells.Item(11, 2) "2024-11-09"
$wsAll.Cells.Item(11, 3).Value = "合肥·一生必听的钢琴曲—“从巴赫 · 莫扎特到肖邦 · 李斯特”钢琴圣手谭小棠独奏音乐会"
$wsAll.Cells.Item(11, 4).Value = "徽州大道辅路与祁门路辅路交叉口北120米 包河凤凰剧院"
$wsAll.Cells.Item(11, 5).Value = "2024.11.09 19:30-11.09 21:00"
$wsAll.Cells.Item(11, 6).Value = 0
$wsAll.Cells.Item(11, 7).Value = 56
$wsAll.Cells.Item(11, 8).Value = "https://show.bilibili.com/platform/detail.html?id=90593"
$wsAll.Cells.Item(11, 9).Value = "//i2.hdslb.com/bfs/openplatform/202408/SYfLxnO21723442234232.jpeg"
